# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (column E, rows 16-30) list gets reversed in order
# (2104 down to 2002 instead of 2002 up to 2104), and the one populated
# date in column F (Salario Basico date helper) moves from the first
# data row to the last one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New descending period list for rows 16..30 (was ascending 2002..2104).
$periodos = @("2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002")

$firstRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

# The lone non-44944 date value in column F moves from row 16 to row 30.
$ws.Cells.Item(16, 6).Value = 38951
$ws.Cells.Item(30, 6).Value = 44944
